$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: global_competitiveness_index.csv / WEF / Global Competitiveness Index
# Change D17 "R script" from wef.R to "not using" (TODO-style comment)
$ws.Range("D17").Value = "not using"

# Row 20: CopyofGPSSpaymentsdata2015draft.xlsx / World Bank / Global Payment Systems Survey
# Add R script "wb.R"
$ws.Range("D20").Value = "wb.R"

# Row 21: Global Findex Database.xlsx / World Bank / Global Findex
# Add R script "wb.R"
$ws.Range("D21").Value = "wb.R"

# Update selection to D22 (matches final selection in diff)
$ws.Range("D22").Select()
